$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 12336307
$ws.Range("I2").Value = 9270535
$ws.Range("J2").Value = 19234292
$ws.Range("K2").Value = 9270535
$ws.Range("L2").Value = 19234292
$ws.Range("M2").Value = -9270422
$ws.Range("N2").Value = -19234518
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("H121").Value = 737.1539
$ws.Range("I121").Value = 329.25
$ws.Range("J121").Value = 811.3182
$ws.Range("K121").Value = 987.75
$ws.Range("L121").Value = 2433.9546
$ws.Range("M121").Value = 759.25
$ws.Range("N121").Value = -5927.9546

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 400
$ws.Range("I3").Value = 400
$ws.Range("K3").Value = 400
$ws.Range("M3").Value = -285
$ws.Range("H17").Value = 2500
$ws.Range("J17").Value = 2500
$ws.Range("L17").Value = 2500
$ws.Range("N17").Value = -2846
$ws.Range("H32").Value = 2161.26
$ws.Range("I32").Value = 1984.809
$ws.Range("J32").Value = 3588.9092
$ws.Range("K32").Value = 1984.809
$ws.Range("L32").Value = 3588.9092
$ws.Range("M32").Value = -1697.809
$ws.Range("N32").Value = -4162.9092
$ws.Range("H61").Value = 2114.182
$ws.Range("I61").Value = 2096.6667
$ws.Range("J61").Value = 2135.2
$ws.Range("K61").Value = 2096.6667
$ws.Range("L61").Value = 2135.2
$ws.Range("M61").Value = -1884.6667
$ws.Range("N61").Value = -2559.2
$ws.Range("H102").Value = 846.41174
$ws.Range("I102").Value = 846.41174
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 846.41174
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 775.58826
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 1018.125
$ws.Range("I122").Value = 982.5625
$ws.Range("J122").Value = 1089.25
$ws.Range("K122").Value = 2947.6875
$ws.Range("L122").Value = 3267.75
$ws.Range("M122").Value = -497.6875
$ws.Range("N122").Value = -8167.75
$ws.Range("H132").Value = 22730902
$ws.Range("I132").Value = 29412808
$ws.Range("J132").Value = 12421.2
$ws.Range("K132").Value = 88238424
$ws.Range("L132").Value = 37263.60000000001
$ws.Range("M132").Value = -88235894
$ws.Range("N132").Value = -42323.60000000001
$ws.Range("H136").Value = 2114.182
$ws.Range("I136").Value = 2096.6667
$ws.Range("J136").Value = 2135.2
$ws.Range("K136").Value = 6290.000100000001
$ws.Range("L136").Value = 6405.599999999999
$ws.Range("M136").Value = -3740.000100000001
$ws.Range("N136").Value = -11505.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 24000
$ws.Range("J16").Value = 24000
$ws.Range("L16").Value = 24000
$ws.Range("N16").Value = -24340
$ws.Range("H99").Value = 1723
$ws.Range("I99").Value = 1102.3077
$ws.Range("J99").Value = 2299.3572
$ws.Range("K99").Value = 1102.3077
$ws.Range("L99").Value = 2299.3572
$ws.Range("M99").Value = 395.6922999999999
$ws.Range("N99").Value = -5295.3572

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 402
$ws.Range("I2").Value = 402
$ws.Range("K2").Value = 402
$ws.Range("M2").Value = -289
$ws.Range("H16").Value = 1971.125
$ws.Range("I16").Value = 950.75
$ws.Range("J16").Value = 2991.5
$ws.Range("K16").Value = 950.75
$ws.Range("L16").Value = 2991.5
$ws.Range("M16").Value = -663.75
$ws.Range("N16").Value = -3565.5
$ws.Range("H113").Value = 1971.125
$ws.Range("I113").Value = 950.75
$ws.Range("J113").Value = 2991.5
$ws.Range("K113").Value = 950.75
$ws.Range("L113").Value = 2991.5
$ws.Range("M113").Value = 1219.25
$ws.Range("N113").Value = -7331.5
$ws.Range("H132").Value = 9012957
$ws.Range("I132").Value = 1011.8571
$ws.Range("J132").Value = 37050120
$ws.Range("K132").Value = 3035.5713
$ws.Range("L132").Value = 111150360
$ws.Range("M132").Value = -505.5712999999996
$ws.Range("N132").Value = -111155420
$ws.Range("H134").Value = 1204.4762
$ws.Range("I134").Value = 1218.8
$ws.Range("J134").Value = 1168.6666
$ws.Range("K134").Value = 3656.4
$ws.Range("L134").Value = 3505.9998
$ws.Range("M134").Value = -1121.4
$ws.Range("N134").Value = -8575.9998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 24825468
$ws.Range("I5").Value = 22222652
$ws.Range("K5").Value = 66667956
$ws.Range("M5").Value = -66667844
$ws.Range("H68").Value = 8052.2144
$ws.Range("I68").Value = 734
$ws.Range("K68").Value = 2202
$ws.Range("M68").Value = -1391
$ws.Range("H71").Value = 8052.2144
$ws.Range("I71").Value = 734
$ws.Range("K71").Value = 6606
$ws.Range("M71").Value = -2550
$ws.Range("H135").Value = 24825468
$ws.Range("I135").Value = 22222652
$ws.Range("K135").Value = 200003868
$ws.Range("M135").Value = -200001333

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 10000
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10224
$ws.Range("H97").Value = 596.1739
$ws.Range("I97").Value = 573.2727
$ws.Range("J97").Value = 1100
$ws.Range("K97").Value = 573.2727
$ws.Range("L97").Value = 1100
$ws.Range("M97").Value = -77.27269999999999
$ws.Range("N97").Value = -2092
$ws.Range("H122").Value = 19237354
$ws.Range("I122").Value = 21745790
$ws.Range("J122").Value = 6000.3335
$ws.Range("K122").Value = 65237370
$ws.Range("L122").Value = 18001.0005
$ws.Range("M122").Value = -65234920
$ws.Range("N122").Value = -22901.0005
$ws.Range("H132").Value = 6637.3477
$ws.Range("I132").Value = 1524.4
$ws.Range("J132").Value = 10570.385
$ws.Range("K132").Value = 4573.200000000001
$ws.Range("L132").Value = 31711.155
$ws.Range("M132").Value = -2043.200000000001
$ws.Range("N132").Value = -36771.155

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1554.5333
$ws.Range("I68").Value = 1603.75
$ws.Range("J68").Value = 1498.2858
$ws.Range("K68").Value = 1603.75
$ws.Range("L68").Value = 1498.2858
$ws.Range("M68").Value = -854.75
$ws.Range("N68").Value = -2996.2858
$ws.Range("H71").Value = 1554.5333
$ws.Range("I71").Value = 1603.75
$ws.Range("J71").Value = 1498.2858
$ws.Range("K71").Value = 8018.75
$ws.Range("L71").Value = 7491.429
$ws.Range("M71").Value = -4274.75
$ws.Range("N71").Value = -14979.429
$ws.Range("H82").Value = 1133.6666
$ws.Range("I82").Value = 1133.6666
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1133.6666
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -772.6666
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 1133.6666
$ws.Range("I85").Value = 1133.6666
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1133.6666
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 114.3334
$ws.Range("N85").ClearContents()
$ws.Range("H100").Value = 3339.9524
$ws.Range("I100").Value = 3126.6667
$ws.Range("J100").Value = 3499.9167
$ws.Range("K100").Value = 3126.6667
$ws.Range("L100").Value = 3499.9167
$ws.Range("M100").Value = -2585.6667
$ws.Range("N100").Value = -4581.9167
$ws.Range("H132").Value = 18524646
$ws.Range("I132").Value = 28573104
$ws.Range("J132").Value = 14326.737
$ws.Range("K132").Value = 85719312
$ws.Range("L132").Value = 42980.211
$ws.Range("M132").Value = -85716782
$ws.Range("N132").Value = -48040.211
